$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 0.350653
$ws.Range("H2").Value = 1.051959
$ws.Range("I2").Value = 0.04536179359243143
$ws.Range("J2").Value = 0.04536179359243143
$ws.Range("M2").Value = 6.045145666666667
$ws.Range("N2").Value = 18.135437
$ws.Range("O2").Value = 0.8160840232643366
$ws.Range("P2").Value = 0.8160840232643367
$ws.Range("Q2").Value = 2.119748463453667
$ws.Range("R2").Value = 19.077736171083
$ws.Range("S2").Value = 0.03701903501739785
$ws.Range("T2").Value = 0.03701903501739785
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 0.350653
$ws.Range("H3").Value = 1.051959
$ws.Range("I3").Value = 0.04536179359243143
$ws.Range("J3").Value = 0.04536179359243143
$ws.Range("O3").Value = 0.09212864864242169
$ws.Range("P3").Value = 0.09212864864242169
$ws.Range("Q3").Value = 0.2393008021633334
$ws.Range("R3").Value = 2.15370721947
$ws.Range("S3").Value = 0.004179120743667171
$ws.Range("T3").Value = 0.004179120743667171
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 0.350653
$ws.Range("H4").Value = 1.051959
$ws.Range("I4").Value = 0.04536179359243143
$ws.Range("J4").Value = 0.04536179359243143
$ws.Range("M4").Value = 0.6799149999999999
$ws.Range("N4").Value = 2.039745
$ws.Range("O4").Value = 0.09178732809324164
$ws.Range("P4").Value = 0.09178732809324165
$ws.Range("Q4").Value = 0.238414234495
$ws.Range("R4").Value = 2.145728110455
$ws.Range("S4").Value = 0.00416363783136641
$ws.Range("T4").Value = 0.004163637831366411
$ws.Range("I5").Value = 0.8482855786262421
$ws.Range("J5").Value = 0.8482855786262421
$ws.Range("M5").Value = 6.045145666666667
$ws.Range("N5").Value = 18.135437
$ws.Range("O5").Value = 0.8160840232643366
$ws.Range("P5").Value = 0.8160840232643367
$ws.Range("Q5").Value = 39.64023265964733
$ws.Range("R5").Value = 356.762093936826
$ws.Range("S5").Value = 0.6922723078824193
$ws.Range("T5").Value = 0.6922723078824194
$ws.Range("I6").Value = 0.8482855786262421
$ws.Range("J6").Value = 0.8482855786262421
$ws.Range("O6").Value = 0.09212864864242169
$ws.Range("P6").Value = 0.09212864864242169
$ws.Range("S6").Value = 0.07815140402169043
$ws.Range("T6").Value = 0.07815140402169043
$ws.Range("I7").Value = 0.8482855786262421
$ws.Range("J7").Value = 0.8482855786262421
$ws.Range("M7").Value = 0.6799149999999999
$ws.Range("N7").Value = 2.039745
$ws.Range("O7").Value = 0.09178732809324164
$ws.Range("P7").Value = 0.09178732809324165
$ws.Range("Q7").Value = 4.458451503889999
$ws.Range("R7").Value = 40.12606353501
$ws.Range("S7").Value = 0.0778618667221322
$ws.Range("T7").Value = 0.07786186672213222
$ws.Range("G8").Value = 0.7457606666666666
$ws.Range("H8").Value = 2.237282
$ws.Range("I8").Value = 0.09647441040198541
$ws.Range("J8").Value = 0.09647441040198541
$ws.Range("M8").Value = 6.045145666666667
$ws.Range("N8").Value = 18.135437
$ws.Range("O8").Value = 0.8160840232643366
$ws.Range("P8").Value = 0.8160840232643367
$ws.Range("Q8").Value = 4.508231862470445
$ws.Range("R8").Value = 40.574086762234
$ws.Range("S8").Value = 0.07873122498290702
$ws.Range("T8").Value = 0.07873122498290704
$ws.Range("G9").Value = 0.7457606666666666
$ws.Range("H9").Value = 2.237282
$ws.Range("I9").Value = 0.09647441040198541
$ws.Range("J9").Value = 0.09647441040198541
$ws.Range("O9").Value = 0.09212864864242169
$ws.Range("P9").Value = 0.09212864864242169
$ws.Range("Q9").Value = 0.5089393952288889
$ws.Range("R9").Value = 4.58045455706
$ws.Range("S9").Value = 0.008888057058909308
$ws.Range("T9").Value = 0.008888057058909308
$ws.Range("G10").Value = 0.7457606666666666
$ws.Range("H10").Value = 2.237282
$ws.Range("I10").Value = 0.09647441040198541
$ws.Range("J10").Value = 0.09647441040198541
$ws.Range("M10").Value = 0.6799149999999999
$ws.Range("N10").Value = 2.039745
$ws.Range("O10").Value = 0.09178732809324164
$ws.Range("P10").Value = 0.09178732809324165
$ws.Range("Q10").Value = 0.5070538636766666
$ws.Range("R10").Value = 4.56348477309
$ws.Range("S10").Value = 0.008855128360169079
$ws.Range("T10").Value = 0.008855128360169081
$ws.Range("G11").Value = 0.07636
$ws.Range("H11").Value = 0.22908
$ws.Range("I11").Value = 0.009878217379341012
$ws.Range("J11").Value = 0.009878217379341012
$ws.Range("M11").Value = 6.045145666666667
$ws.Range("N11").Value = 18.135437
$ws.Range("O11").Value = 0.8160840232643366
$ws.Range("P11").Value = 0.8160840232643367
$ws.Range("Q11").Value = 0.4616073231066667
$ws.Range("R11").Value = 4.154465907960001
$ws.Range("S11").Value = 0.008061455381612305
$ws.Range("T11").Value = 0.008061455381612305
$ws.Range("G12").Value = 0.07636
$ws.Range("H12").Value = 0.22908
$ws.Range("I12").Value = 0.009878217379341012
$ws.Range("J12").Value = 0.009878217379341012
$ws.Range("O12").Value = 0.09212864864242169
$ws.Range("P12").Value = 0.09212864864242169
$ws.Range("Q12").Value = 0.05211137293333334
$ws.Range("R12").Value = 0.4690023564
$ws.Range("S12").Value = 0.0009100668181547716
$ws.Range("T12").Value = 0.0009100668181547716
$ws.Range("G13").Value = 0.07636
$ws.Range("H13").Value = 0.22908
$ws.Range("I13").Value = 0.009878217379341012
$ws.Range("J13").Value = 0.009878217379341012
$ws.Range("M13").Value = 0.6799149999999999
$ws.Range("N13").Value = 2.039745
$ws.Range("O13").Value = 0.09178732809324164
$ws.Range("P13").Value = 0.09178732809324165
$ws.Range("Q13").Value = 0.05191830939999999
$ws.Range("R13").Value = 0.4672647846
$ws.Range("S13").Value = 0.000906695179573935
$ws.Range("T13").Value = 0.0009066951795739352

Write-Output "Applied all cell updates"
